$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceVolume {
    param(
        [int]$Row,
        $Price,
        $ForceText,
        $Volume
    )
    if ($Price -ne $null) {
        $cell = $ws.Range("D$Row")
        if ($ForceText -eq $true) {
            $cell.NumberFormat = "@"
            $cell.Value = $Price
            $cell.Style = "Normal"
        } else {
            $cell.Value = $Price
        }
    }
    if ($Volume -ne $null) {
        $ws.Range("E$Row").Value = "  $Volume  "
    }
}

Set-PriceVolume 2 "60.705.88" $false "-3.59%"
Set-PriceVolume 3 "2.904.97" $false "-4.13%"
Set-PriceVolume 4 $null $false "-0.05%"
Set-PriceVolume 5 "587.97" $true "-1.02%"
Set-PriceVolume 6 "143.97" $true "-6.33%"
Set-PriceVolume 8 "0.502" $true "-2.37%"
Set-PriceVolume 9 "2.903.79" $false "-4.13%"
Set-PriceVolume 10 "6.68" $true "-3.33%"
Set-PriceVolume 11 $null $false "-5.00%"
Set-PriceVolume 12 "0.442" $true "-4.44%"
Set-PriceVolume 13 $null $false "-3.82%"
Set-PriceVolume 14 "33.38" $true "-6.14%"
Set-PriceVolume 15 $null $false "+1.61%"
Set-PriceVolume 16 "3.386.74" $false "-4.30%"
Set-PriceVolume 17 "60.689.11" $false "-3.59%"
Set-PriceVolume 18 "6.67" $true "-5.81%"
Set-PriceVolume 19 "2.906.18" $false "-4.10%"
Set-PriceVolume 20 "427.69" $true "-5.58%"
Set-PriceVolume 21 "13.50" $true "-5.53%"
Set-PriceVolume 22 "0.683" $true "-2.32%"
Set-PriceVolume 23 "7.06" $true "-6.21%"
Set-PriceVolume 24 "80.81" $true "-2.81%"
Set-PriceVolume 25 "10.79" $true "-5.77%"
Set-PriceVolume 26 $null $false "-5.45%"
Set-PriceVolume 27 "11.86" $true "-4.63%"
Set-PriceVolume 28 $null $false "+0.06%"
Set-PriceVolume 29 $null $false "-0.03%"
Set-PriceVolume 30 $null $false "-3.59%"
Set-PriceVolume 31 "7.16" $true "-4.83%"
Set-PriceVolume 32 $null $false "-3.73%"
Set-PriceVolume 33 "26.38" $true "-4.47%"
Set-PriceVolume 34 $null $false "-3.96%"
Set-PriceVolume 35 "0.0₃0856" $false "-0.70%"
Set-PriceVolume 36 $null $false "-3.21%"
Set-PriceVolume 37 $null $false "-5.70%"
Set-PriceVolume 38 $null $false "-4.04%"
Set-PriceVolume 39 "49.38" $true "-2.17%"
Set-PriceVolume 40 $null $false "-6.14%"
Set-PriceVolume 41 $null $false "-5.93%"
Set-PriceVolume 42 $null $false "-5.91%"
Set-PriceVolume 43 "0.294" $true "-5.15%"
Set-PriceVolume 44 "41.17" $true "-7.77%"
Set-PriceVolume 45 "0.0349" $true "-3.04%"
Set-PriceVolume 46 "372.68" $true "-5.46%"
Set-PriceVolume 47 "2.697.18" $false "-0.88%"
Set-PriceVolume 48 "132.19" $true "-0.48%"
Set-PriceVolume 50 $null $false "-5.74%"
Set-PriceVolume 51 $null $false "-3.06%"
